# Applies:
#  - Sheet "Informes_Gestió_Versions_2023-09-11_2023-09-14": append a new
#    data row (row 33) describing the "Pegats" (patches) deployment mail.
#  - Sheet "Resumen": rename the "Otros" row to "Pegats", add a brand new
#    "Otros" row right after it, move the "Total" row down one row with
#    refreshed totals, and update the narrative sentence below the table.
#  - Chart "Resum": extend the three series' category/value ranges to
#    include the new "Resumen" row.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Informes_Gestió_Versions_2023-09-11_2023-09-14")
$ws2 = $wb.Worksheets.Item("Resumen")

# ---------------------------------------------------------------------
# 1) New deployment row in the detail sheet (row 33), cloning row 31's
#    look (unshaded "odd" row with the date column formatted as a date).
# ---------------------------------------------------------------------
$ws1.Range("A31:I31").Copy($ws1.Range("A33:I33"))
$excel.CutCopyMode = 0

$ws1.Range("A33").Value = "PRODUCCION"
$ws1.Range("B33").Value = "Mail final. RFC C1642769: Distribucio pegats seguretat anual - 2023"
$ws1.Range("C33").Value = "Pegats"
$ws1.Range("D33").Value = "OK"
$ws1.Range("E33").Value = "NO"
$ws1.Range("I33").Value = 45183

# ---------------------------------------------------------------------
# 2) Resumen table: shift the "Total" row down to row 10 (carrying its
#    shaded style), insert a fresh "Otros" row at row 9, and relabel the
#    existing row 8 from "Otros" to "Pegats".
# ---------------------------------------------------------------------
$ws2.Range("A8:E8").Copy()
$ws2.Range("A10:E10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws2.Range("A10").Value = "Total"
$ws2.Range("B10").Value = 27
$ws2.Range("C10").Value = 5
$ws2.Range("D10").Value = 32
$ws2.Range("E10").Value = 9

$ws2.Range("A9").Value = "Otros"
$ws2.Range("B9").Value = 1
$ws2.Range("C9").Value = 0
$ws2.Range("D9").Value = 1
$ws2.Range("E9").Value = 0

$ws2.Range("A8").Value = "Pegats"

$ws2.Range("A17").Value = "De 32 desplegaments, 9 han sigut urgents, d" + [char]8217 + "aquests, 0 tenien incidència associada."

# ---------------------------------------------------------------------
# 3) Chart series now span through the new Total row (row 10).
# ---------------------------------------------------------------------
$co = $ws2.ChartObjects().Item(1)
$chart = $co.Chart
$chart.SeriesCollection().Item(1).Formula = "=SERIES(Resumen!`$B`$1,Resumen!`$A`$2:`$A`$10,Resumen!`$B`$2:`$B`$10,1)"
$chart.SeriesCollection().Item(2).Formula = "=SERIES(Resumen!`$C`$1,Resumen!`$A`$2:`$A`$10,Resumen!`$C`$2:`$C`$10,2)"
$chart.SeriesCollection().Item(3).Formula = "=SERIES(Resumen!`$D`$1,Resumen!`$A`$2:`$A`$10,Resumen!`$D`$2:`$D`$10,3)"
